$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Climate")

# Touch the pre-existing "0.000"-formatted cell (D3) first with the new
# 4-decimal format so the engine reuses its existing style slot in place
# once the whole B2:D6 block below is restyled to match it.
$ws.Range("D3").NumberFormat = "0.0000"

# Use newest weather inputs (monthly means), including rain.rate means.
$ws.Cells.Item(2, 2).Value = 4.2751882725795
$ws.Cells.Item(2, 3).Value = 4.04186842718138
$ws.Cells.Item(2, 4).Value = 0.0556576719496239

$ws.Cells.Item(3, 2).Value = 8.19748979894988
$ws.Cells.Item(3, 3).Value = 3.81357212945869
$ws.Cells.Item(3, 4).Value = 0.116938350042087

$ws.Cells.Item(4, 2).Value = 12.3751389561464
$ws.Cells.Item(4, 3).Value = 3.45567498501387
$ws.Cells.Item(4, 4).Value = 0.071780395514249

$ws.Cells.Item(5, 2).Value = 16.8185301623819
$ws.Cells.Item(5, 3).Value = 3.10883875876392
$ws.Cells.Item(5, 4).Value = 0.107050334741296

$ws.Cells.Item(6, 2).Value = 14.3971459030677
$ws.Cells.Item(6, 3).Value = 3.30879945612945
$ws.Cells.Item(6, 4).Value = 0.128927644374806

# All climate figures now share the same 4-decimal-place display format.
$ws.Range("B2:D6").NumberFormat = "0.0000"

# Move the active selection on the Climate sheet.
$ws.Activate()
$ws.Range("C10").Select()
